# Add a new row of sample data to the "Samples" sheet describing the
# second spectramax data file, then leave the selection where the
# author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Samples")

$ws.Range("A3").Value = "Plate reader"
$ws.Range("B3").Value = "`$GITHUB_WORKSPACE/test/inputs/spectramax-data2.txt"
$ws.Range("C3").Value = "600,700,530_485,530_485_1,530_485_2"

$ws.Range("E16").Select()
